# Update "想去人数" (F column) values on sheet "展览" (Worksheets(1))
# and sheet "全部类型" (Worksheets(4)) to reflect refreshed counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14979
$ws1.Range("F3").Value = 18915
$ws1.Range("F5").Value = 135
$ws1.Range("F14").Value = 146
$ws1.Range("F22").Value = 7872
$ws1.Range("F26").Value = 61
$ws1.Range("F27").Value = 1237
$ws1.Range("F29").Value = 6026
$ws1.Range("F30").Value = 111
$ws1.Range("F35").Value = 5392
$ws1.Range("F37").Value = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14979
$ws4.Range("F3").Value = 18915
$ws4.Range("F5").Value = 135
$ws4.Range("F14").Value = 146
$ws4.Range("F23").Value = 7872
$ws4.Range("F27").Value = 61
$ws4.Range("F28").Value = 1237
$ws4.Range("F32").Value = 6026
$ws4.Range("F33").Value = 111
$ws4.Range("F38").Value = 5392
$ws4.Range("F40").Value = 6
